# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G holds the "K" values; update the rows whose K value changed after
# the regeneration of std/mean and recalculation of s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kUpdates = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 0
    8  = 3
    9  = 2
    10 = 2
    12 = 2
    13 = 1
    14 = 1
    16 = 2
    17 = 1
}

foreach ($row in $kUpdates.Keys) {
    $ws.Range("G$row").Value = $kUpdates[$row]
}
